# Atualizado por script em 02-12-2023 20:45
#
# This script:
#   1. Swaps the match content (columns F:V) between row pairs 16/17, 24/25
#      and 31/32 - the two fixtures in each pair had been recorded in the
#      wrong order, so we exchange everything except the shared index/date
#      columns (A:E stay where they are).
#   2. Appends four newly-scraped fixtures as rows 38-41 (index 37-40),
#      copying the formatting of an existing data row and then filling in
#      the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param($ws, $row1, $row2)

    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range($col + $row1).Value()
        $vals2[$col] = $ws.Range($col + $row2).Value()
    }

    foreach ($col in $cols) {
        $ws.Range($col + $row1).Value = $vals2[$col]
        $ws.Range($col + $row2).Value = $vals1[$col]
    }
}

function Set-MatchRow {
    param($ws, $rowNum, $values)

    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

    # Bring over the same cell formatting used by the other data rows
    # (bold/centered/bordered index in col A, datetime format in col E).
    $ws.Range("A2:V2").Copy()
    $ws.Range("A" + $rowNum + ":V" + $rowNum).PasteSpecial(-4122)

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowNum).Value = $values[$i]
    }
}

# --- 1. Fix the swapped fixtures -------------------------------------------------
Swap-MatchRows $ws 16 17
Swap-MatchRows $ws 24 25
Swap-MatchRows $ws 31 32

# --- 2. Append the newly scraped fixtures ---------------------------------------
Set-MatchRow $ws 38 @(37, "israel", "ligat-ha-al", "2023-2024", 45262.65625, "Hapoel Hadera", 3, "Maccabi Petah Tikva", 0, 3.02, "05/10/2023 02:42", 3.94, "02/12/2023 15:06", 3.38, "05/10/2023 02:42", 3.41, "02/12/2023 15:23", 2.25, "05/10/2023 02:42", 2, "02/12/2023 15:06", "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-hadera-maccabi-petah-tikva/UNlPKRI8/")
Set-MatchRow $ws 39 @(38, "israel", "ligat-ha-al", "2023-2024", 45262.6875, "Sakhnin", 1, "Netanya", 1, 3.96, "01/10/2023 19:42", 3.7, "02/12/2023 11:00", 3.49, "01/10/2023 19:42", 3.37, "02/12/2023 11:00", 1.93, "01/10/2023 19:42", 2.09, "02/12/2023 11:00", "https://www.betexplorer.com/football/israel/ligat-ha-al/sakhnin-netanya/tOq5BNue/")
Set-MatchRow $ws 40 @(39, "israel", "ligat-ha-al", "2023-2024", 45262.75, "Maccabi Bnei Raina", 0, "SC Ashdod", 0, 2.21, "30/09/2023 19:12", 2.12, "02/12/2023 17:51", 3.31, "30/09/2023 19:12", 3.41, "02/12/2023 17:51", 3.33, "30/09/2023 19:12", 3.57, "02/12/2023 17:51", "https://www.betexplorer.com/football/israel/ligat-ha-al/maccabi-bnei-raina-sc-ashdod/AZ5MJoYE/")
Set-MatchRow $ws 41 @(40, "israel", "ligat-ha-al", "2023-2024", 45262.77083333334, "Hapoel Petah Tikva", 1, "Beitar Jerusalem", 0, 3.21, "04/10/2023 16:12", 3.15, "02/12/2023 18:29", 3.36, "04/10/2023 16:12", 3.38, "02/12/2023 18:29", 2.25, "04/10/2023 16:12", 2.32, "02/12/2023 15:48", "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-petah-tikva-beitar-jerusalem/rm4UHP2R/")
